$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.236.53'
$ws.Range("E2").Value = '  +3.63%  '

# Row 3
$ws.Range("D3").Value = '2.280.48'
$ws.Range("E3").Value = '  +2.73%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.24'
$ws.Range("E5").Value = '  +1.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.03'
$ws.Range("E6").Value = '  +5.53%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.591'
$ws.Range("E7").Value = '  +3.20%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("E9").Value = '  +3.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.29'
$ws.Range("E10").Value = '  +4.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +3.23%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.92'
$ws.Range("E12").Value = '  +3.49%  '

# Row 13
$ws.Range("E13").Value = '  +3.03%  '

# Row 14
$ws.Range("D14").Value = '2.623.39'
$ws.Range("E14").Value = '  +2.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.880'
$ws.Range("E15").Value = '  +2.99%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.62'
$ws.Range("E16").Value = '  +4.26%  '

# Row 17
$ws.Range("D17").Value = '2.282.94'
$ws.Range("E17").Value = '  +5.34%  '

# Row 18
$ws.Range("D18").Value = '44.114.79'
$ws.Range("E18").Value = '  +3.73%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.65'
$ws.Range("E19").Value = '  +7.21%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000100'
$ws.Range("E20").Value = '  +4.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  +3.16%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.43'
$ws.Range("E22").Value = '  +1.30%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.23'
$ws.Range("E23").Value = '  +1.41%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '239.99'
$ws.Range("E24").Value = '  +1.82%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  +5.45%  '

# Row 26
$ws.Range("E26").Value = '  -0.25%  '

# Row 27
$ws.Range("E27").Value = '  +2.76%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.26'
$ws.Range("E28").Value = '  +1.45%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.59'
$ws.Range("E29").Value = '  +17.11%  '

# Row 30
$ws.Range("E30").Value = '  +1.58%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.54'
$ws.Range("E31").Value = '  +2.00%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.64'
$ws.Range("E32").Value = '  +0.84%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0889'
$ws.Range("E33").Value = '  +1.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '161.59'
$ws.Range("E34").Value = '  +2.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.42'
$ws.Range("E35").Value = '  +6.39%  '

# Row 36
$ws.Range("E36").Value = '  +0.81%  '

# Row 37
$ws.Range("E37").Value = '  +9.23%  '

# Row 38
$ws.Range("E38").Value = '  +0.42%  '

# Row 39
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.99'
$ws.Range("E39").Value = '  +13.77%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.53'
$ws.Range("E40").Value = '  +2.17%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.108'
$ws.Range("E41").Value = '  +5.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.04'
$ws.Range("E42").Value = '  +33.93%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0330'
$ws.Range("E43").Value = '  +3.72%  '

# Row 44
$ws.Range("E44").Value = '  +0.07%  '

# Row 45
$ws.Range("D45").Value = '1.817.37'
$ws.Range("E45").Value = '  +1.99%  '

# Row 46
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.87'
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.210'
$ws.Range("E47").Value = '  +2.16%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.45'
$ws.Range("E48").Value = '  +2.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.90'
$ws.Range("E49").Value = '  -1.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.88'
$ws.Range("E50").Value = '  +5.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.18'
$ws.Range("E51").Value = '  +0.01%  '
